# Optima_template.xlsx edit - "Regeneration of en/fr/es country data books
# after child program paras updated"
#
# Substantive change: updated input data (coverage/effectiveness fractions)
# for the "Programmes pour les enfants" sheet, rows 2-53 (columns D-H), plus
# the resulting recalculation of the dependent shared-formula rows below
# them (57-163) which the engine recomputes automatically.
#
# Cosmetic changes replicated where the COM surface allows it: which sheet
# tab is selected/active (moves from "Donnees pop de l'annee de ref" to
# "Dependances du programme"), and the selected range on
# "Programmes pour les enfants" (now D2:H53).

$wb = $excel.ActiveWorkbook

$wsChildren = $wb.Worksheets.Item("Programmes pour les enfants")

$wsChildren.Range("F2").Value = 0.39473684210526322
$wsChildren.Range("G2").Value = 0.39473684210526322
$wsChildren.Range("H2").Value = 0.39473684210526322
$wsChildren.Range("F3").Value = 0.30769230769230765
$wsChildren.Range("G3").Value = 0.30769230769230765
$wsChildren.Range("H3").Value = 0.30769230769230765
$wsChildren.Range("F4").Value = 0.38507462686567184
$wsChildren.Range("G4").Value = 0.38507462686567184
$wsChildren.Range("H4").Value = 0.38507462686567184
$wsChildren.Range("F6").Value = 0.25970149253731345
$wsChildren.Range("G6").Value = 0.25970149253731345
$wsChildren.Range("F8").Value = 0.25970149253731345
$wsChildren.Range("G8").Value = 0.25970149253731345
$wsChildren.Range("F10").Value = 0.25970149253731345
$wsChildren.Range("G10").Value = 0.25970149253731345
$wsChildren.Range("F12").Value = 0.25970149253731345
$wsChildren.Range("G12").Value = 0.25970149253731345
$wsChildren.Range("F14").Value = 0.25970149253731345
$wsChildren.Range("G14").Value = 0.25970149253731345
$wsChildren.Range("F16").Value = 0.25970149253731345
$wsChildren.Range("G16").Value = 0.25970149253731345
$wsChildren.Range("F18").Value = 0.7
$wsChildren.Range("F20").Value = 0.84
$wsChildren.Range("D21").Value = 0.28260869565217389
$wsChildren.Range("F21").Value = 0
$wsChildren.Range("F22").Value = 0
$wsChildren.Range("D23").Value = 0.28260869565217389
$wsChildren.Range("F23").Value = 0
$wsChildren.Range("F24").Value = 0
$wsChildren.Range("D25").Value = 0.28260869565217389
$wsChildren.Range("F25").Value = 0
$wsChildren.Range("F26").Value = 0
$wsChildren.Range("F27").Value = 1
$wsChildren.Range("F28").Value = 0
$wsChildren.Range("F29").Value = 0
$wsChildren.Range("F30").Value = 1
$wsChildren.Range("F31").Value = 0
$wsChildren.Range("F32").Value = 0
$wsChildren.Range("F33").Value = 1
$wsChildren.Range("F34").Value = 0
$wsChildren.Range("F35").Value = 0
$wsChildren.Range("F36").Value = 1
$wsChildren.Range("F37").Value = 0
$wsChildren.Range("F38").Value = 0
$wsChildren.Range("F39").Value = 1
$wsChildren.Range("F40").Value = 0
$wsChildren.Range("F41").Value = 0
$wsChildren.Range("F42").Value = 0.3
$wsChildren.Range("F43").Value = 0.5
$wsChildren.Range("F44").Value = 0.65
$wsChildren.Range("F45").Value = 0.3
$wsChildren.Range("F46").Value = 0.49
$wsChildren.Range("F47").Value = 0.52
$wsChildren.Range("F48").Value = 0.88
$wsChildren.Range("D49").Value = 0.78409090909090906
$wsChildren.Range("E49").Value = 0.78409090909090906
$wsChildren.Range("F49").Value = 0.78409090909090906
$wsChildren.Range("G49").Value = 0.78409090909090906
$wsChildren.Range("H49").Value = 0.78409090909090906
$wsChildren.Range("D50").Value = 0.88372093023255816
$wsChildren.Range("E50").Value = 0.88372093023255816
$wsChildren.Range("F50").Value = 0.88372093023255816
$wsChildren.Range("G50").Value = 0.88372093023255816
$wsChildren.Range("H50").Value = 0.88372093023255816
$wsChildren.Range("F51").Value = 0.86
$wsChildren.Range("F52").Value = 0
$wsChildren.Range("F53").Value = 0

# --- Window / selection state -------------------------------------------
# Move the active range on the child-programmes sheet to D2:H53 (the block
# of inputs that was just updated), then hand tab-selection over to
# "Dependances du programme" (sheet index 10, 0-based) as in the target
# workbook, leaving "Donnees pop de l'annee de ref" no longer the selected
# tab.

$wsChildren.Activate()
$wsChildren.Range("D2:H53").Select()

$wsDependencies = $wb.Worksheets.Item("Dépendances du programme")
$wsDependencies.Select()
